# Update "想去人数" (number of people interested) figures on the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet.
#
# These two sheets list overlapping events; the "全部类型" sheet has one
# extra row (an 演出/concert entry) inserted at row 8, so the same
# logical rows live at F4:F12 on 展览 and F4:F7 + F9:F13 on 全部类型.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 247
$ws1.Range("F5").Value  = 2934
$ws1.Range("F6").Value  = 1997
$ws1.Range("F7").Value  = 385
$ws1.Range("F8").Value  = 134
$ws1.Range("F9").Value  = 1084
$ws1.Range("F10").Value = 203
$ws1.Range("F11").Value = 498
$ws1.Range("F12").Value = 55

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 247
$ws4.Range("F5").Value  = 2934
$ws4.Range("F6").Value  = 1997
$ws4.Range("F7").Value  = 385
$ws4.Range("F9").Value  = 134
$ws4.Range("F10").Value = 1084
$ws4.Range("F11").Value = 203
$ws4.Range("F12").Value = 498
$ws4.Range("F13").Value = 55
